$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "25-03-2025"
$ws.Range("B6").Value = "Gujarat Titans vs Punjab Kings"
$ws.Range("C6").Value = "Gujarat Titans"
$ws.Range("D6").Value = "Gujarat Titans"
